# The deck currently uses the "Integral" (Red Violet) theme on its
# slide master / presentation theme part, while the notes master's
# theme part still carries the stock "Office Theme" colour scheme.
#
# The authored edit swaps the two: the presentation's active theme
# becomes the standard "Office Theme" colour scheme (the one that used
# to live only on the notes master), while the previously-active
# "Integral" / "Red Violet" colours move off to the side.
#
# Apply that by rewriting the 12 theme colour slots on the
# presentation's (only) slide master / design to the stock Office
# theme RGB values, in the canonical clrScheme slot order:
#   dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# dk1 = 000000
$colors.Item(1).RGB = 0
# lt1 = FFFFFF
$colors.Item(2).RGB = 16777215
# dk2 = 44546A
$colors.Item(3).RGB = 6968388
# lt2 = E7E6E6
$colors.Item(4).RGB = 15132391
# accent1 = 5B9BD5
$colors.Item(5).RGB = 13998939
# accent2 = ED7D31
$colors.Item(6).RGB = 3243501
# accent3 = A5A5A5
$colors.Item(7).RGB = 10855845
# accent4 = FFC000
$colors.Item(8).RGB = 49407
# accent5 = 4472C4
$colors.Item(9).RGB = 12874308
# accent6 = 70AD47
$colors.Item(10).RGB = 4697456
# hlink = 0563C1
$colors.Item(11).RGB = 12673797
# folHlink = 954F72
$colors.Item(12).RGB = 7491477
